# Refreshed IMF WEO data with APR 25 forecasts
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet ---
$ws.Name = "WEO_Data_Denmark"

# --- Updated source / footnote text (shared strings) ---
$ws.Range("E2").Value = "Source: National Statistics Office Latest actual data: 2023 National accounts manual used: European System of Accounts (ESA) 2010 GDP valuation: Market prices Reporting in calendar year: Yes Start/end months of reporting year: January/December Base year: 2020 Chain-weighted: Yes, from 1980 Primary domestic currency: Danish krone Data last updated: 04/09/2025"

$ws.Range("E3").Value = "Source: National Statistics Office Latest actual data: 2024 Harmonized prices: Yes Base year: 2015 Primary domestic currency: Danish krone Data last updated: 04/09/2025"
$ws.Range("E5").Value = "Source: National Statistics Office Latest actual data: 2024 Harmonized prices: Yes Base year: 2015 Primary domestic currency: Danish krone Data last updated: 04/09/2025"

$ws.Range("E7").Value = "Source: National Statistics Office Latest actual data: 2024 Primary domestic currency: Danish krone Data last updated: 04/09/2025"

$ws.Range("A9").Value = "International Monetary Fund, World Economic Outlook Database, April 2025"

# --- Row 1: year headers. BD1 becomes the 2030 year column, notes marker shifts to BE1 ---
$ws.Range("BD1").Value = 2030
$ws.Range("BE1").Value = "Estimates Start After"

# --- Row 2: Gross domestic product, current prices (National currency, Billions) ---
$ws.Range("AX2").Value = 2960.8870000000002
$ws.Range("AY2").Value = 3111.2629999999999
$ws.Range("AZ2").Value = 3230.97
$ws.Range("BA2").Value = 3349.7
$ws.Range("BB2").Value = 3464.9749999999999
$ws.Range("BC2").Value = 3586.9050000000002
$ws.Range("BD2").NumberFormat = "#,##0.00"
$ws.Range("BD2").Value = 3710.3530000000001
$ws.Range("BE2").Value = 2023

# --- Row 3: Inflation, average consumer prices (Index) ---
$ws.Range("AX3").Value = 119.133
$ws.Range("AY3").Value = 121.43899999999999
$ws.Range("AZ3").Value = 123.99
$ws.Range("BA3").Value = 126.46899999999999
$ws.Range("BB3").Value = 128.999
$ws.Range("BC3").Value = 131.57900000000001
$ws.Range("BD3").Value = 134.21
$ws.Range("BE3").Value = 2024

# --- Row 4: Inflation, average consumer prices (Percent change) ---
$ws.Range("AX4").Value = 1.268
$ws.Range("AY4").Value = 1.9359999999999999
$ws.Range("AZ4").Value = 2.1
$ws.Range("BD4").Value = 2
$ws.Range("BE4").Value = 2024

# --- Row 5: Inflation, end of period consumer prices (Index) ---
$ws.Range("AW5").Value = 118
$ws.Range("AX5").Value = 120.2
$ws.Range("AY5").Value = 122.744
$ws.Range("AZ5").Value = 125.322
$ws.Range("BA5").Value = 127.828
$ws.Range("BB5").Value = 130.38399999999999
$ws.Range("BC5").Value = 132.99199999999999
$ws.Range("BD5").Value = 135.65199999999999
$ws.Range("BE5").Value = 2024

# --- Row 6: Inflation, end of period consumer prices (Percent change) ---
$ws.Range("AW6").Value = 0.42599999999999999
$ws.Range("AX6").Value = 1.8640000000000001
$ws.Range("AY6").Value = 2.1160000000000001
$ws.Range("AZ6").Value = 2.1
$ws.Range("BD6").Value = 2
$ws.Range("BE6").Value = 2024

# --- Row 7: Population (Persons, Millions) ---
$ws.Range("AX7").Value = 5.9610000000000003
$ws.Range("AY7").Value = 6.0019999999999998
$ws.Range("AZ7").Value = 6.0460000000000003
$ws.Range("BA7").Value = 6.093
$ws.Range("BB7").Value = 6.15
$ws.Range("BC7").Value = 6.1840000000000002
$ws.Range("BD7").Value = 6.2270000000000003
$ws.Range("BE7").Value = 2024
